$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2's data into a new row 5, but change the "Experience (Years)" value (G) to 8
$ws.Range("A5").Value2 = $ws.Range("A2").Value2
$ws.Range("C5").Value2 = $ws.Range("C2").Value2
$ws.Range("D5").Value2 = $ws.Range("D2").Value2
$ws.Range("E5").Value2 = $ws.Range("E2").Value2
$ws.Range("F5").Value2 = $ws.Range("F2").Value2
$ws.Range("G5").Value2 = 8
$ws.Range("H5").Value2 = $ws.Range("H2").Value2

# Give column D (Email) an explicit custom width, as recorded in the diff
$ws.Columns.Item(4).ColumnWidth = 16.67

# Select the newly active cell to match the workbook's recorded selection state
$ws.Range("G5").Select()
